$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column D: "canonical SMILES" ---

# Copy formatting of the header cell (C2) onto the new header cell (D2)
$ws.Range("C2").Copy() | Out-Null
$ws.Range("D2").PasteSpecial(-4122) | Out-Null
$ws.Range("D2").Value = "canonical SMILES"

# For each data row, copy C's formatting to D, then mirror C's SMILES value into D
for ($r = 3; $r -le 10; $r++) {
    $srcCell = "C$r"
    $dstCell = "D$r"
    $ws.Range($srcCell).Copy() | Out-Null
    $ws.Range($dstCell).PasteSpecial(-4122) | Out-Null
    $ws.Range($dstCell).Value = $ws.Range($srcCell).Value2
}

$excel.CutCopyMode = 0

# Set the width for the new column (closest achievable value to the source's 36.85546875)
$ws.Columns.Item(4).ColumnWidth = 36
